# Auto update Excel log
$wb = $excel.ActiveWorkbook

function Set-TextValue($cell, $value) {
    # Force the cell to keep the literal text instead of letting Excel
    # auto-convert look-alike dates/percentages into numeric values.
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

# --- PIR sheet: append rows 42-59 ---
$pirRows = @(
    @("2026-02-01", "13:49:31", "13:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-02-01", "13:49:32", "13:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-02-01", "13:49:34", "13:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-02-01", "13:49:34", "13:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-02-01", "13:49:34", "13:00", "Bathroom", "Motion Detected", "Active"),
    @("2026-02-01", "13:49:34", "13:00", "Bathroom", "Motion Detected", "Active"),
    @("2026-02-01", "13:50:54", "13:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-02-01", "13:50:59", "13:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-02-01", "13:51:04", "13:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-02-01", "13:51:09", "13:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-02-01", "13:51:14", "13:00", "Bathroom", "Motion Detected", "Active"),
    @("2026-02-01", "13:51:22", "13:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-02-01", "13:51:27", "13:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-02-01", "13:51:32", "13:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-02-01", "13:51:37", "13:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-02-01", "13:51:42", "13:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-02-01", "13:51:47", "13:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-02-01", "13:51:52", "13:00", "Bathroom", "No Motion", "Inactive")
)
$ws = $wb.Worksheets.Item("PIR")
$startRow = 42
for ($i = 0; $i -lt $pirRows.Count; $i++) {
    $row = $pirRows[$i]
    $r = $startRow + $i
    Set-TextValue $ws.Cells.Item($r, 1) $row[0]
    for ($c = 1; $c -lt 6; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $row[$c]
    }
}

# --- Humidity sheet: append rows 18-30 ---
# Column A looks like a date ("2026-02-01") and column E looks like a
# percentage ("81.4%"); force both to text so Excel keeps the literal
# strings instead of converting them into numeric values.
$humidityRows = @(
    @("2026-02-01", "13:49:31", "13:00", "Bathroom", "81.4%", "Active"),
    @("2026-02-01", "13:49:31", "13:00", "Bathroom", "81.2%", "Active"),
    @("2026-02-01", "13:49:33", "13:00", "Bathroom", "80.2%", "Active"),
    @("2026-02-01", "13:49:34", "13:00", "Bathroom", "81.0%", "Active"),
    @("2026-02-01", "13:50:58", "13:00", "Bathroom", "80.4%", "Active"),
    @("2026-02-01", "13:51:03", "13:00", "Bathroom", "79.4%", "Active"),
    @("2026-02-01", "13:51:08", "13:00", "Bathroom", "80.3%", "Active"),
    @("2026-02-01", "13:51:13", "13:00", "Bathroom", "79.4%", "Active"),
    @("2026-02-01", "13:51:18", "13:00", "Bathroom", "80.3%", "Active"),
    @("2026-02-01", "13:51:23", "13:00", "Bathroom", "79.6%", "Active"),
    @("2026-02-01", "13:51:28", "13:00", "Bathroom", "80.4%", "Active"),
    @("2026-02-01", "13:51:38", "13:00", "Bathroom", "79.5%", "Active"),
    @("2026-02-01", "13:51:48", "13:00", "Bathroom", "79.4%", "Active")
)
$ws = $wb.Worksheets.Item("Humidity")
$startRow = 18
for ($i = 0; $i -lt $humidityRows.Count; $i++) {
    $row = $humidityRows[$i]
    $r = $startRow + $i
    Set-TextValue $ws.Cells.Item($r, 1) $row[0]
    for ($c = 1; $c -lt 6; $c++) {
        if ($c -eq 4) {
            Set-TextValue $ws.Cells.Item($r, $c + 1) $row[$c]
        } else {
            $ws.Cells.Item($r, $c + 1).Value = $row[$c]
        }
    }
}

# --- Proximity sheet: append row 22 ---
$ws = $wb.Worksheets.Item("Proximity")
$r = 22
Set-TextValue $ws.Cells.Item($r, 1) "2026-02-01"
$ws.Cells.Item($r, 2).Value = "13:49:32"
$ws.Cells.Item($r, 3).Value = "13:00"
$ws.Cells.Item($r, 4).Value = "Living Room Main Door"
$ws.Cells.Item($r, 5).Value = "ENTER"
$ws.Cells.Item($r, 6).Value = "User ENTERED Living Room Main Door"
